$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in accuracy (C) / f1 (D) values for the SVM block (rows 24-32) ---
$ws.Range("C24").Value = 0.26
$ws.Range("D24").Value = 0.26

$ws.Range("C25").Value = 0.36
$ws.Range("D25").Value = 0.37

$ws.Range("C26").Value = 0.42
$ws.Range("D26").Value = 0.42

$ws.Range("C27").Value = 0.3
$ws.Range("D27").Value = 0.3

$ws.Range("C28").Value = 0.4
$ws.Range("D28").Value = 0.4

$ws.Range("C29").Value = 0.43
$ws.Range("D29").Value = 0.43

$ws.Range("C30").Value = 0.2
$ws.Range("D30").Value = 0.22

$ws.Range("C31").Value = 0.44
$ws.Range("D31").Value = 0.44

$ws.Range("C32").Value = 0.33
$ws.Range("D32").Value = 0.33

# --- Add new "3 layers nn" results block header (row 39) ---
$ws.Range("A39").Value = "3 layers nn"
$ws.Range("B39").Value = "parameter"
$ws.Range("C39").Value = "accuracy"
$ws.Range("D39").Value = "f1"
$ws.Range("E39").Value = "cm"

# Match formatting used by the other section headers (row 1 / A23 / A33):
# bold font, thin border, left-aligned text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A39:E39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Update the active selection to reflect the new bottom of the sheet ---
$ws.Range("A40").Select() | Out-Null
